$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 540.4
$ws.Range("I19").Value = 483.16666
$ws.Range("J19").Value = 626.25
$ws.Range("K19").Value = 483.16666
$ws.Range("L19").Value = 626.25
$ws.Range("M19").Value = -308.16666
$ws.Range("N19").Value = -976.25
$ws.Range("H98").Value = 2913.9092
$ws.Range("I98").Value = 2077.353
$ws.Range("J98").Value = 5758.2
$ws.Range("K98").Value = 2077.353
$ws.Range("L98").Value = 5758.2
$ws.Range("M98").Value = -579.3530000000001
$ws.Range("N98").Value = -8754.200000000001
$ws.Range("H116").Value = 6677.909
$ws.Range("I116").Value = 10084.583
$ws.Range("J116").Value = 2589.9
$ws.Range("K116").Value = 10084.583
$ws.Range("L116").Value = 2589.9
$ws.Range("M116").Value = -6642.583000000001
$ws.Range("N116").Value = -9473.9
$ws.Range("H122").Value = 2913.9092
$ws.Range("I122").Value = 2077.353
$ws.Range("J122").Value = 5758.2
$ws.Range("K122").Value = 6232.059
$ws.Range("L122").Value = 17274.6
$ws.Range("M122").Value = -3782.059
$ws.Range("N122").Value = -22174.6
$ws.Range("H132").Value = 1332.3143
$ws.Range("I132").Value = 1107.7587
$ws.Range("J132").Value = 2417.6667
$ws.Range("K132").Value = 3323.2761
$ws.Range("L132").Value = 7253.000100000001
$ws.Range("M132").Value = -793.2761
$ws.Range("N132").Value = -12313.0001
$ws.Range("H137").Value = 1451.36
$ws.Range("I137").Value = 1250.5555
$ws.Range("J137").Value = 1967.7142
$ws.Range("K137").Value = 3751.6665
$ws.Range("L137").Value = 5903.142599999999
$ws.Range("M137").Value = -1201.6665
$ws.Range("N137").Value = -11003.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7380.3687
$ws.Range("I61").Value = 7380.3687
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7380.3687
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -7168.3687
$ws.Range("H63").Value = 142859840
$ws.Range("I63").Value = 142859840
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 142859840
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -142859154
$ws.Range("H66").Value = 142859840
$ws.Range("I66").Value = 142859840
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 714299200
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -714295768
$ws.Range("H136").Value = 7380.3687
$ws.Range("I136").Value = 7380.3687
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 22141.1061
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -19591.1061

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6460.393
$ws.Range("I31").Value = 1735.8235
$ws.Range("J31").Value = 13762
$ws.Range("K31").Value = 1735.8235
$ws.Range("L31").Value = 13762
$ws.Range("M31").Value = -1440.8235
$ws.Range("N31").Value = -14352
$ws.Range("H34").Value = 6460.393
$ws.Range("I34").Value = 1735.8235
$ws.Range("J34").Value = 13762
$ws.Range("K34").Value = 1735.8235
$ws.Range("L34").Value = 13762
$ws.Range("M34").Value = -1533.8235
$ws.Range("N34").Value = -14166
$ws.Range("H58").Value = 1357.2745
$ws.Range("I58").Value = 624.37933
$ws.Range("J58").Value = 2323.3635
$ws.Range("K58").Value = 624.37933
$ws.Range("L58").Value = 2323.3635
$ws.Range("M58").Value = -421.37933
$ws.Range("N58").Value = -2729.3635
$ws.Range("H87").Value = 23000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 23000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 23000
$ws.Range("N87").Value = -25372
$ws.Range("H90").Value = 23000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 23000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 69000
$ws.Range("N90").Value = -80856
$ws.Range("H136").Value = 1357.2745
$ws.Range("I136").Value = 624.37933
$ws.Range("J136").Value = 2323.3635
$ws.Range("K136").Value = 1873.13799
$ws.Range("L136").Value = 6970.0905
$ws.Range("M136").Value = 676.8620100000001
$ws.Range("N136").Value = -12070.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 14.8
$ws.Range("I2").Value = 16.5
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 99
$ws.Range("L2").Value = 48
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = -274
$ws.Range("H132").Value = 1795389.8
$ws.Range("I132").Value = 933.3333
$ws.Range("J132").Value = 1886633.2
$ws.Range("K132").Value = 8399.9997
$ws.Range("L132").Value = 16979698.8
$ws.Range("M132").Value = -5869.9997
$ws.Range("N132").Value = -16984758.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5818.5
$ws.Range("I126").Value = 7398.8887
$ws.Range("J126").Value = 2262.625
$ws.Range("K126").Value = 22196.6661
$ws.Range("L126").Value = 6787.875
$ws.Range("M126").Value = -19726.6661
$ws.Range("N126").Value = -11727.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 111115544
$ws.Range("I40").Value = 125004110
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 125004110
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -125003974
$ws.Range("N40").Value = -7272
$ws.Range("H132").Value = 45846576
$ws.Range("I132").Value = 45846576
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 137539728
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -137537198

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H122").Value = 1624.875
$ws.Range("I122").Value = 1435.6364
$ws.Range("J122").Value = 1785
$ws.Range("K122").Value = 4306.9092
$ws.Range("L122").Value = 5355
$ws.Range("M122").Value = -1856.9092
$ws.Range("N122").Value = -10255
$ws.Range("H126").Value = 987.7778
$ws.Range("I126").Value = 955.7143
$ws.Range("J126").Value = 1100
$ws.Range("K126").Value = 2867.1429
$ws.Range("L126").Value = 3300
$ws.Range("M126").Value = -397.1428999999998
$ws.Range("N126").Value = -8240
$ws.Range("H132").Value = 2465.077
$ws.Range("I132").Value = 1620.4
$ws.Range("J132").Value = 2993
$ws.Range("K132").Value = 4861.200000000001
$ws.Range("L132").Value = 8979
$ws.Range("M132").Value = -2331.200000000001
$ws.Range("N132").Value = -14039
$ws.Range("H136").Value = 3132.8572
$ws.Range("I136").Value = 4584
$ws.Range("J136").Value = 1875.2
$ws.Range("K136").Value = 13752
$ws.Range("L136").Value = 5625.6
$ws.Range("M136").Value = -11202
$ws.Range("N136").Value = -10725.6
